$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E5 run 9 sample data, dated 2021-07-08 (diff adds sheet row 45).
$ws.Range("A45").Value = 20210708
$ws.Range("B45").Value = 2228.9570833859798
$ws.Range("C45").Value = 2224.4699999999998
$ws.Range("D45").Formula = "=100*(B45-C45)/C45"
$ws.Range("E45").Value = 180
$ws.Range("F45").Value = "CRM OPENED 20210526"

# The old per-batch "CRM OPENED 20210624"/"CRM OPENED 20210706" notes on
# rows 43/44 are replaced with the recurring "CRM OPENED 20210526" note
# (those two strings are no longer used anywhere in the sheet).
$ws.Range("F43").Value = "CRM OPENED 20210526"
$ws.Range("F44").Value = "CRM OPENED 20210526"

# Column C (Batch value) gets an explicit custom width now that it holds a
# wider note column to its right being referenced more / sheet widened.
$ws.Columns("C").ColumnWidth = 14.5

# Leave the cursor on the newly-entered formula cell, matching the saved
# selection state.
$ws.Range("D45").Select()
